$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block (columns K:P), entered in the order that makes new
# shared strings land at the same indices as the target workbook:
# 10 "prog (lots of tapping)", 11 "Raspberry Pi zero",
# 12 "Idel or Program running?", 13 "Realtime", 14 "buffer overload".

# Row 7 first (introduces "prog (lots of tapping)" then later "buffer overload")
$ws.Range("K7").Value = 0.069
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 0
$ws.Range("N7").Formula = "=(K7*3600)/(L7*60 +M7)"
$ws.Range("O7").Value = "prog (lots of tapping)"

# Row 1 header labels
$ws.Range("K1").Value = "Raspberry Pi zero"
$ws.Range("O1").Value = "Idel or Program running?"
$ws.Range("P1").Value = "Realtime"

# Back to row 7 for the last new string
$ws.Range("P7").Value = "buffer overload"

# Row 4 sub-headers (reuse existing shared strings Wh/m/s)
$ws.Range("K4").Value = "Wh"
$ws.Range("L4").Value = "m"
$ws.Range("M4").Value = "s"

# Row 5
$ws.Range("K5").Value = 0.05
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 3
$ws.Range("N5").Formula = "=(K5*3600)/(L5*60 +M5)"
$ws.Range("O5").Value = "idle"

# Row 6
$ws.Range("K6").Value = 0.1438
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 45
$ws.Range("N6").Formula = "=(K6*3600)/(L6*60 +M6)"
$ws.Range("O6").Value = "prog (only detector)"

# --- Convert E14:E21 into a single shared formula, matching the target ---
$ws.Range("E14:E21").Formula = "=(B14*3600)/(C14*60 +D14)"

# --- Column widths for the new best-fit columns ---
$ws.Columns("K").ColumnWidth = 14.166666666666666
$ws.Columns("N").ColumnWidth = 11.333333333333334
$ws.Columns("O").ColumnWidth = 19.666666666666668
$ws.Columns("P").ColumnWidth = 12.333333333333334

# --- Selection cursor moves to R9 in the saved file ---
$ws.Range("R9").Select()
